# Update "Use Cases:" to "Use Cases & Salient Features :" in the
# "Use Cases / Business Value" summary box on slide 3.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$para = $tr.Paragraphs(1, 1)
$run = $para.Runs(1, 1)
$run.Text = "Use Cases & Salient Features :"
